# Fixed university images not showing.
# Adds a new "image" column to the University sheet's Table1, and fills in
# the university-code abbreviation used to look up each institution's logo.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("University")

# --- 1. Extend Table1 with a 6th column ("image") ---------------------
$lo = $ws.ListObjects.Item("Table1")
$null = $lo.ListColumns.Add()

# --- 2. Populate the new column ----------------------------------------
# Header first, then each first-occurrence of a code (so that new shared
# strings end up created in the same order as the source workbook), then
# the remaining repeated values.
$ws.Range("F1").Value = "image"

$ws.Range("F138").Value = "KMITL"
$ws.Range("F143").Value = "NU"
$ws.Range("F116").Value = "KU"
$ws.Range("F191").Value = "RMU"
$ws.Range("F153").Value = "RU"
$ws.Range("F125").Value = "MU"
$ws.Range("F127").Value = "SU"
$ws.Range("F128").Value = "PSU"
$ws.Range("F123").Value = "TU"
$ws.Range("F117").Value = "KKU"
$ws.Range("F118").Value = "CU"
$ws.Range("F119").Value = "CMU"
$ws.Range("F120").Value = "KMUTT"
$ws.Range("F121").Value = "KMUTNB"
$ws.Range("F122").Value = "SUT"
$ws.Range("F124").Value = "BUU"
$ws.Range("F126").Value = "SWU"

$ws.Range("F154").Value = "RU"
$ws.Range("F155").Value = "RU"
$ws.Range("F156").Value = "RU"
$ws.Range("F157").Value = "RU"
$ws.Range("F158").Value = "RU"
$ws.Range("F159").Value = "RU"
$ws.Range("F160").Value = "RU"
$ws.Range("F161").Value = "RU"
$ws.Range("F162").Value = "RU"
$ws.Range("F163").Value = "RU"
$ws.Range("F164").Value = "RU"
$ws.Range("F165").Value = "RU"
$ws.Range("F166").Value = "RU"
$ws.Range("F167").Value = "RU"
$ws.Range("F168").Value = "RU"
$ws.Range("F169").Value = "RU"
$ws.Range("F170").Value = "RU"
$ws.Range("F171").Value = "RU"
$ws.Range("F172").Value = "RU"
$ws.Range("F173").Value = "RU"
$ws.Range("F174").Value = "RU"
$ws.Range("F175").Value = "RU"
$ws.Range("F176").Value = "RU"
$ws.Range("F177").Value = "RU"
$ws.Range("F178").Value = "RU"
$ws.Range("F179").Value = "RU"
$ws.Range("F180").Value = "RU"
$ws.Range("F181").Value = "RU"
$ws.Range("F182").Value = "RU"
$ws.Range("F183").Value = "RU"
$ws.Range("F184").Value = "RU"
$ws.Range("F185").Value = "RU"
$ws.Range("F186").Value = "RU"
$ws.Range("F187").Value = "RU"
$ws.Range("F188").Value = "RU"
$ws.Range("F189").Value = "RU"
$ws.Range("F190").Value = "RU"

$ws.Range("F192").Value = "RMU"
$ws.Range("F193").Value = "RMU"
$ws.Range("F194").Value = "RMU"
$ws.Range("F195").Value = "RMU"
$ws.Range("F196").Value = "RMU"
$ws.Range("F197").Value = "RMU"
$ws.Range("F198").Value = "RMU"
$ws.Range("F199").Value = "RMU"

# --- 3. Restore the author's on-screen selection/scroll position -------
$ws.Activate()
$ws.Range("F143").Select()
